$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 12555
$ws.Range("F3").Value = 7026
$ws.Range("F6").Value = 441
$ws.Range("F10").Value = 978
$ws.Range("F11").Value = 130
$ws.Range("F12").Value = 337
$ws.Range("F13").Value = 988
$ws.Range("F17").Value = 229
$ws.Range("F18").Value = 357
$ws.Range("F19").Value = 16
$ws.Range("F20").Value = 266
$ws.Range("F21").Value = 296
$ws.Range("F23").Value = 115
$ws.Range("F24").Value = 357
$ws.Range("F25").Value = 5172
$ws.Range("F26").Value = 65
$ws.Range("F27").Value = 1382
$ws.Range("F28").Value = 292
$ws.Range("F29").Value = 1218
$ws.Range("F30").Value = 1316
$ws.Range("F31").Value = 582
$ws.Range("F33").Value = 3714

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 92
$ws.Range("F7").Value = 36
$ws.Range("F16").Value = 18

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9217
$ws.Range("F3").Value = 546
$ws.Range("F4").Value = 1956

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9217
$ws.Range("F3").Value = 546
$ws.Range("F4").Value = 1956
$ws.Range("F5").Value = 12555
$ws.Range("F6").Value = 7026
$ws.Range("F7").Value = 92
$ws.Range("F10").Value = 441
$ws.Range("F14").Value = 978
$ws.Range("F15").Value = 130
$ws.Range("F16").Value = 337
$ws.Range("F17").Value = 988
$ws.Range("F21").Value = 229
$ws.Range("F22").Value = 357
$ws.Range("F23").Value = 16
$ws.Range("F24").Value = 266
$ws.Range("F25").Value = 296
$ws.Range("F31").Value = 357
$ws.Range("F32").Value = 5172
$ws.Range("F33").Value = 65
$ws.Range("F34").Value = 1382
$ws.Range("F37").Value = 292
$ws.Range("F39").Value = 1218
$ws.Range("F40").Value = 1316
$ws.Range("F41").Value = 582
$ws.Range("F46").Value = 3714
$ws.Range("F47").Value = 18
